$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = 0.8303227957310438
$ws.Range("J3").Value = 0.5326835134425084
$ws.Range("K3").Value = 0.6985624488779159
$ws.Range("L3").Value = 2.96954242281555
